$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 96.666664
$ws.Range("K2").Value = 96.666664
$ws.Range("M2").Value = 16.333336
$ws.Range("H5").Value = 229.85715
$ws.Range("I5").Value = 229.85715
$ws.Range("K5").Value = 229.85715
$ws.Range("M5").Value = -114.85715
$ws.Range("H9").Value = 291.42856
$ws.Range("I9").Value = 275
$ws.Range("K9").Value = 275
$ws.Range("M9").Value = -106
$ws.Range("H40").Value = 7226.8184
$ws.Range("I40").Value = 3749.5
$ws.Range("K40").Value = 3749.5
$ws.Range("M40").Value = -3574.5
$ws.Range("H103").Value = 781.625
$ws.Range("I103").Value = 770.6
$ws.Range("K103").Value = 2311.8
$ws.Range("M103").Value = -1725.8
$ws.Range("H127").Value = 492.15
$ws.Range("I127").Value = 386.4737
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 1159.4211
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = 3800.5789
$ws.Range("N127").Value = -17420
$ws.Range("H129").Value = 2594.25
$ws.Range("I129").Value = 1792.3334
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 5377.0002
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -377.0002000000004
$ws.Range("N129").Value = -25000
$ws.Range("H137").Value = 2127.5
$ws.Range("I137").Value = 2153
$ws.Range("K137").Value = 6459
$ws.Range("M137").Value = -3909

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 487.75
$ws.Range("I4").Value = 448.66666
$ws.Range("K4").Value = 448.66666
$ws.Range("M4").Value = -332.66666
$ws.Range("H5").Value = 100
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("N5").Value = -324
$ws.Range("H32").Value = 2793.182
$ws.Range("I32").Value = 2872.5
$ws.Range("K32").Value = 2872.5
$ws.Range("M32").Value = -2585.5
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -330
$ws.Range("H22").Value = 399.875
$ws.Range("I22").Value = 399.875
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 399.875
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -226.875
$ws.Range("N22").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1989
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1989
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 5967
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10867
$ws.Range("H132").Value = 2610.5454
$ws.Range("I132").Value = 2301.9
$ws.Range("K132").Value = 6905.700000000001
$ws.Range("M132").Value = -4375.700000000001
$ws.Range("H134").Value = 2498
$ws.Range("I134").Value = 2498
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7494
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4959
$ws.Range("N134").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 312805.94
$ws.Range("I4").Value = 312805.94
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 938417.8200000001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -938305.8200000001
$ws.Range("N4").ClearContents()
$ws.Range("H11").Value = 8335769
$ws.Range("I11").Value = 10002888
$ws.Range("K11").Value = 30008664
$ws.Range("M11").Value = -30008524
$ws.Range("H68").Value = 2235.3333
$ws.Range("I68").Value = 1200
$ws.Range("K68").Value = 3600
$ws.Range("M68").Value = -2789
$ws.Range("H71").Value = 2235.3333
$ws.Range("I71").Value = 1200
$ws.Range("K71").Value = 10800
$ws.Range("M71").Value = -6744

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1783.8462
$ws.Range("I102").Value = 1653.6364
$ws.Range("K102").Value = 1653.6364
$ws.Range("M102").Value = -31.63640000000009

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8800
$ws.Range("J2").Value = 8800
$ws.Range("L2").Value = 8800
$ws.Range("N2").Value = -9024
$ws.Range("H7").Value = 5500
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H46").Value = 3531.5625
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 3700.3333
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 3700.3333
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -4076.3333
$ws.Range("H68").Value = 2780
$ws.Range("I68").Value = 2780
$ws.Range("K68").Value = 2780
$ws.Range("M68").Value = -2031
$ws.Range("H71").Value = 2780
$ws.Range("I71").Value = 2780
$ws.Range("K71").Value = 13900
$ws.Range("M71").Value = -10156
$ws.Range("H82").Value = 2762.5
$ws.Range("I82").Value = 1975
$ws.Range("J82").Value = 3550
$ws.Range("K82").Value = 1975
$ws.Range("L82").Value = 3550
$ws.Range("M82").Value = -1614
$ws.Range("N82").Value = -4272
$ws.Range("H85").Value = 2762.5
$ws.Range("I85").Value = 1975
$ws.Range("J85").Value = 3550
$ws.Range("K85").Value = 1975
$ws.Range("L85").Value = 3550
$ws.Range("M85").Value = -727
$ws.Range("N85").Value = -6046
$ws.Range("H122").Value = 9499.5
$ws.Range("I122").Value = 9499.5
$ws.Range("K122").Value = 28498.5
$ws.Range("M122").Value = -26048.5
$ws.Range("H123").Value = 99995
$ws.Range("J123").Value = 99995
$ws.Range("L123").Value = 99995
$ws.Range("N123").Value = -109795
$ws.Range("H126").Value = 5500
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 27967.166
$ws.Range("I81").Value = 31560.6
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 63121.2
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -62060.2
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 27967.166
$ws.Range("I84").Value = 31560.6
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 315606
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -310302
$ws.Range("N84").Value = -110608
$ws.Range("H124").Value = 50429
$ws.Range("J124").Value = 50429
$ws.Range("L124").Value = 50429
$ws.Range("N124").Value = -60249
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
